# Weekly update: insert a new week's worth of price records (3 rows, one per
# quality grade: Primera/Segunda/Tercera) at the top of the data block
# (row 966), pushing the existing rows down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 966..968 (shifts old 966.. down to 969..)
$ws.Range("A966:A968").EntireRow.Insert()

# --- Row 966 : Primera ---
$ws.Range("A966").Value = 6
$ws.Range("B966").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C966").Value = "Metropolitana"
$ws.Range("D966").Value = 44783
$ws.Range("E966").Value = 13
$ws.Range("F966").Value = 100114014
$ws.Range("G966").Value = "Betarraga"
$ws.Range("H966").Value = "Sin especificar"
$ws.Range("I966").Value = "Primera"
$ws.Range("J966").Value = 17000
$ws.Range("K966").Value = 170
$ws.Range("L966").Value = 170
$ws.Range("M966").Value = 170
$ws.Range("N966").Value = "`$/unidad"
$ws.Range("O966").Value = "Región Metropolitana"
$ws.Range("P966").Value = 170
$ws.Range("Q966").Value = 1
$ws.Range("R966").Value = "Hortaliza"

# --- Row 967 : Segunda ---
$ws.Range("A967").Value = 6
$ws.Range("B967").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C967").Value = "Metropolitana"
$ws.Range("D967").Value = 44783
$ws.Range("E967").Value = 13
$ws.Range("F967").Value = 100114014
$ws.Range("G967").Value = "Betarraga"
$ws.Range("H967").Value = "Sin especificar"
$ws.Range("I967").Value = "Segunda"
$ws.Range("J967").Value = 14000
$ws.Range("K967").Value = 120
$ws.Range("L967").Value = 120
$ws.Range("M967").Value = 120
$ws.Range("N967").Value = "`$/unidad"
$ws.Range("O967").Value = "Región Metropolitana"
$ws.Range("P967").Value = 120
$ws.Range("Q967").Value = 1
$ws.Range("R967").Value = "Hortaliza"

# --- Row 968 : Tercera ---
$ws.Range("A968").Value = 6
$ws.Range("B968").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C968").Value = "Metropolitana"
$ws.Range("D968").Value = 44783
$ws.Range("E968").Value = 13
$ws.Range("F968").Value = 100114014
$ws.Range("G968").Value = "Betarraga"
$ws.Range("H968").Value = "Sin especificar"
$ws.Range("I968").Value = "Tercera"
$ws.Range("J968").Value = 11000
$ws.Range("K968").Value = 85
$ws.Range("L968").Value = 85
$ws.Range("M968").Value = 85
$ws.Range("N968").Value = "`$/unidad"
$ws.Range("O968").Value = "Región Metropolitana"
$ws.Range("P968").Value = 85
$ws.Range("Q968").Value = 1
$ws.Range("R968").Value = "Hortaliza"
